$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in event name: QSL -> SQL Saturday Lima 2024 (1096)
$ws.Range("B31").Value = "SQL Saturday Lima 2024 (1096)"

# Fill in attendance numbers for SQL Saturday Oregon 2024 (1082) row 29
$ws.Range("C29").Value = 331
$ws.Range("D29").Value = 240

# Nudge the shared "no show rate" formula so it recalculates against the
# newly entered values (re-applying the same formula forces a fresh calc).
$ws.Range("E29").Formula = $ws.Range("E29").Formula

# Update the active selection cell to D29
$ws.Range("D29").Select()
